# Replace the long templated letter body with a short "Hi there!" greeting.
# The target paragraph currently holds a single run made up of alternating
# <w:br/> line breaks and <w:t/> text pieces:
#
#   <br/><br/>Dear [Employer],<br/><br/>...<br/><br/>[Doctor Name]
#
# We need to keep the two leading line breaks but collapse everything from
# "Dear [Employer]," through "[Doctor Name]" into the single sentence
# "Hi there!".

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()

$found = $find.Execute(
    "Dear \[Employer\],*\[Doctor Name\]",  # wildcard pattern spanning the whole block
    $false,    # MatchCase
    $false,    # MatchWholeWord
    $true,     # MatchWildcards
    $false,    # MatchSoundsLike
    $false,    # MatchAllWordForms
    $true,     # Forward
    1,         # Wrap (wdFindContinue)
    $false,    # Format
    "Hi there!",  # ReplaceWith
    2          # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find the expected letter body to replace."
}
